# "Added functionality to Search Bar" — appends newly-searched/added records
# to the Stores sheet and their corresponding (still-blank) placeholder rows
# on the Electric sheet.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($ws, [int]$row, [int]$col, [string]$text)

    # Force the cell to store a literal text value (matches the source
    # workbook, where every column is written as text regardless of
    # whether it "looks" numeric). A leading apostrophe is Excel's
    # standard quote-prefix idiom for this; plain alphabetic text does
    # not need it, so we only add it when required, keeping unrelated
    # cells on the default style.
    $needsPrefix = $false
    if ($text -eq "") {
        $needsPrefix = $true
    } else {
        $needsPrefix = $text -match '^[+-]?[0-9]+(\.[0-9]+)?$'
    }

    if ($needsPrefix) {
        $ws.Cells.Item($row, $col).Value = "'" + $text
    } else {
        $ws.Cells.Item($row, $col).Value = $text
    }
}

# ---------------------------------------------------------------------
# Stores sheet: append newly searched / added store records (rows 8-18)
# ---------------------------------------------------------------------
$storesWs = $wb.Worksheets.Item("Stores")

$storeRows = @(
    @("7",  "1ST FLOOR", "54",   "4"),
    @("8",  "1ST FLOOR", "32",   "42"),
    @("9",  "2ND FLOOR", "234",  "24"),
    @("10", "CHICKEN",   "23",   "3"),
    @("11", "1ST FLOOR", "5",    "3"),
    @("12", "VEGETABLE", "3",    "wew"),
    @("13", "1ST FLOOR", "asd",  "asd"),
    @("14", "FRUIT",     "23",   "24"),
    @("15", "GROCERY",   "34",   "24"),
    @("16", "1ST FLOOR", "wer",  "wer"),
    @("17", "SPECIAL",   "Spec", "sp")
)

$r = 8
foreach ($rowVals in $storeRows) {
    $col = 1
    foreach ($val in $rowVals) {
        Set-TextValue $storesWs $r $col $val
        $col = $col + 1
    }
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Electric sheet: extend the (still blank) data rows down to row 16
# ---------------------------------------------------------------------
$electricWs = $wb.Worksheets.Item("Electric")

$row = 14
while ($row -le 16) {
    $col = 1
    while ($col -le 6) {
        Set-TextValue $electricWs $row $col ""
        $col = $col + 1
    }
    $row = $row + 1
}
